$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 264.125
$ws.Range("I11").Value = 264.125
$ws.Range("K11").Value = 264.125
$ws.Range("M11").Value = -124.125
$ws.Range("H33").Value = 725
$ws.Range("I33").Value = 725
$ws.Range("K33").Value = 725
$ws.Range("M33").Value = -496
$ws.Range("H43").Value = 2408.0908
$ws.Range("I43").Value = 2633.3333
$ws.Range("J43").Value = 2137.8
$ws.Range("K43").Value = 2633.3333
$ws.Range("L43").Value = 2137.8
$ws.Range("M43").Value = -2564.3333
$ws.Range("N43").Value = -2275.8
$ws.Range("H51").Value = 5744.875
$ws.Range("I51").Value = 4094.5
$ws.Range("J51").Value = 5894.909
$ws.Range("K51").Value = 4094.5
$ws.Range("L51").Value = 5894.909
$ws.Range("M51").Value = -3610.5
$ws.Range("N51").Value = -6862.909
$ws.Range("H76").Value = 6861.4614
$ws.Range("J76").Value = 9358.571
$ws.Range("L76").Value = 9358.571
$ws.Range("N76").Value = -9988.571
$ws.Range("H79").Value = 6861.4614
$ws.Range("J79").Value = 9358.571
$ws.Range("L79").Value = 9358.571
$ws.Range("N79").Value = -11542.571
$ws.Range("H88").Value = 2756.182
$ws.Range("J88").Value = 2501.8572
$ws.Range("L88").Value = 2501.8572
$ws.Range("N88").Value = -3313.8572
$ws.Range("H91").Value = 2756.182
$ws.Range("J91").Value = 2501.8572
$ws.Range("L91").Value = 2501.8572
$ws.Range("N91").Value = -5309.8572
$ws.Range("H132").Value = 2943.4
$ws.Range("I132").Value = 1850.356
$ws.Range("J132").Value = 13691.667
$ws.Range("K132").Value = 5551.068
$ws.Range("L132").Value = 41075.001
$ws.Range("M132").Value = -3021.068
$ws.Range("N132").Value = -46135.001
$ws.Range("H137").Value = 2404.5833
$ws.Range("I137").Value = 2222.1428
$ws.Range("K137").Value = 6666.428400000001
$ws.Range("M137").Value = -4116.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2021.65
$ws.Range("I74").Value = 1932.2
$ws.Range("K74").Value = 1932.2
$ws.Range("M74").Value = -1058.2
$ws.Range("H77").Value = 2021.65
$ws.Range("I77").Value = 1932.2
$ws.Range("K77").Value = 9661
$ws.Range("M77").Value = -5293
$ws.Range("H97").Value = 2313.3157
$ws.Range("I97").Value = 1568.0714
$ws.Range("J97").Value = 4400
$ws.Range("K97").Value = 1568.0714
$ws.Range("L97").Value = 4400
$ws.Range("M97").Value = -1072.0714
$ws.Range("N97").Value = -5392
$ws.Range("H122").Value = 2667.3
$ws.Range("I122").Value = 2758.0264
$ws.Range("K122").Value = 8274.0792
$ws.Range("M122").Value = -5824.0792
$ws.Range("H132").Value = 6496287.5
$ws.Range("I132").Value = 3010.5625
$ws.Range("J132").Value = 38463188
$ws.Range("K132").Value = 9031.6875
$ws.Range("L132").Value = 115389564
$ws.Range("M132").Value = -6501.6875
$ws.Range("N132").Value = -115394624

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 79271
$ws.Range("J63").Value = 79271
$ws.Range("L63").Value = 79271
$ws.Range("N63").Value = -80643
$ws.Range("H66").Value = 79271
$ws.Range("J66").Value = 79271
$ws.Range("L66").Value = 237813
$ws.Range("N66").Value = -244677
$ws.Range("H69").Value = 77999
$ws.Range("J69").Value = 99999
$ws.Range("L69").Value = 99999
$ws.Range("N69").Value = -101497
$ws.Range("H72").Value = 77999
$ws.Range("J72").Value = 99999
$ws.Range("L72").Value = 299997
$ws.Range("N72").Value = -307485
$ws.Range("H86").Value = 16850.777
$ws.Range("I86").Value = 19950
$ws.Range("K86").Value = 19950
$ws.Range("M86").Value = -18827
$ws.Range("H89").Value = 16850.777
$ws.Range("I89").Value = 19950
$ws.Range("K89").Value = 99750
$ws.Range("M89").Value = -94134
$ws.Range("H132").Value = 2801.8262
$ws.Range("I132").Value = 2342.1892
$ws.Range("K132").Value = 7026.567599999999
$ws.Range("M132").Value = -4496.567599999999
$ws.Range("H134").Value = 2657.0962
$ws.Range("I134").Value = 2670.6592
$ws.Range("J134").Value = 2582.5
$ws.Range("K134").Value = 8011.9776
$ws.Range("L134").Value = 7747.5
$ws.Range("M134").Value = -5476.9776
$ws.Range("N134").Value = -12817.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 768.6
$ws.Range("I5").Value = 788.2941
$ws.Range("J5").Value = 742.8461
$ws.Range("K5").Value = 2364.8823
$ws.Range("L5").Value = 2228.5383
$ws.Range("M5").Value = -2252.8823
$ws.Range("N5").Value = -2452.5383
$ws.Range("H17").Value = 6344.2
$ws.Range("I17").Value = 60
$ws.Range("K17").Value = 180
$ws.Range("M17").Value = -11
$ws.Range("H18").Value = 487.8
$ws.Range("I18").Value = 430.8889
$ws.Range("K18").Value = 1292.6667
$ws.Range("M18").Value = -1123.6667
$ws.Range("H46").Value = 3737.9
$ws.Range("I46").Value = 499
$ws.Range("K46").Value = 1497
$ws.Range("M46").Value = -1406
$ws.Range("H122").Value = 41048.89
$ws.Range("J122").Value = 7540
$ws.Range("L122").Value = 67860
$ws.Range("N122").Value = -72760
$ws.Range("H129").Value = 1091376.6
$ws.Range("I129").Value = 2361.6667
$ws.Range("J129").Value = 2279392.8
$ws.Range("K129").Value = 7085.000100000001
$ws.Range("L129").Value = 6838178.399999999
$ws.Range("M129").Value = -2085.000100000001
$ws.Range("N129").Value = -6848178.399999999
$ws.Range("H132").Value = 1000
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H135").Value = 768.6
$ws.Range("I135").Value = 788.2941
$ws.Range("J135").Value = 742.8461
$ws.Range("K135").Value = 7094.6469
$ws.Range("L135").Value = 6685.6149
$ws.Range("M135").Value = -4559.6469
$ws.Range("N135").Value = -11755.6149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7802.8184
$ws.Range("I70").Value = 6674.1
$ws.Range("K70").Value = 6674.1
$ws.Range("M70").Value = -6404.1
$ws.Range("H73").Value = 7802.8184
$ws.Range("I73").Value = 6674.1
$ws.Range("K73").Value = 6674.1
$ws.Range("M73").Value = -5738.1
$ws.Range("H132").Value = 2670.3215
$ws.Range("I132").Value = 2631.5334
$ws.Range("J132").Value = 2715.077
$ws.Range("K132").Value = 7894.600199999999
$ws.Range("L132").Value = 8145.231000000001
$ws.Range("M132").Value = -5364.600199999999
$ws.Range("N132").Value = -13205.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2749.5
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 1149.125
$ws.Range("I22").Value = 1118.6
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 1118.6
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -823.5999999999999
$ws.Range("N22").Value = -1790
$ws.Range("H27").Value = 1149.125
$ws.Range("I27").Value = 1118.6
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 1118.6
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -1011.6
$ws.Range("N27").Value = -1414
$ws.Range("H55").Value = 1348.2
$ws.Range("I55").Value = 1447.6
$ws.Range("J55").Value = 1298.5
$ws.Range("K55").Value = 1447.6
$ws.Range("L55").Value = 1298.5
$ws.Range("M55").Value = -1274.6
$ws.Range("N55").Value = -1644.5
$ws.Range("H132").Value = 3511.5833
$ws.Range("I132").Value = 3155.1875
$ws.Range("K132").Value = 9465.5625
$ws.Range("M132").Value = -6935.5625
$ws.Range("H136").Value = 4307.636
$ws.Range("I136").Value = 4373.75
$ws.Range("K136").Value = 13121.25
$ws.Range("M136").Value = -10571.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 7331.3335
$ws.Range("I37").Value = 5997.5
$ws.Range("J37").Value = 9999
$ws.Range("K37").Value = 5997.5
$ws.Range("L37").Value = 9999
$ws.Range("M37").Value = -5794.5
$ws.Range("N37").Value = -10405
$ws.Range("H81").Value = 1674.091
$ws.Range("I81").Value = 1341.5
$ws.Range("K81").Value = 2683
$ws.Range("M81").Value = -1622
$ws.Range("H84").Value = 1674.091
$ws.Range("I84").Value = 1341.5
$ws.Range("K84").Value = 13415
$ws.Range("M84").Value = -8111
$ws.Range("H132").Value = 2428.2444
$ws.Range("I132").Value = 1917.9062
$ws.Range("K132").Value = 5753.7186
$ws.Range("M132").Value = -3223.7186
